$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Egf"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.06368266666666667
$ws.Range("H2").Value = 0.191048
$ws.Range("I2").Value = 0.07976548992723123
$ws.Range("J2").Value = 0.07976548992723123
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.179771666666667
$ws.Range("N2").Value = 6.539315
$ws.Range("O2").Value = 0.2349306639444428
$ws.Range("P2").Value = 0.2349306639444428
$ws.Range("Q2").Value = 0.1388136724577778
$ws.Range("R2").Value = 1.24932305212
$ws.Range("S2").Value = 0.0187393595084582
$ws.Range("T2").Value = 0.01873935950845819

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Egf"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.06368266666666667
$ws.Range("H3").Value = 0.191048
$ws.Range("I3").Value = 0.07976548992723123
$ws.Range("J3").Value = 0.07976548992723123
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.237801
$ws.Range("N3").Value = 12.713403
$ws.Range("O3").Value = 0.4567402255103586
$ws.Range("P3").Value = 0.4567402255103586
$ws.Range("Q3").Value = 0.2698744684826667
$ws.Range("R3").Value = 2.428870216344
$ws.Range("S3").Value = 0.03643210785730783
$ws.Range("T3").Value = 0.03643210785730783

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Egf"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.06368266666666667
$ws.Range("H4").Value = 0.191048
$ws.Range("I4").Value = 0.07976548992723123
$ws.Range("J4").Value = 0.07976548992723123
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.860789
$ws.Range("N4").Value = 8.582367000000001
$ws.Range("O4").Value = 0.3083291105451987
$ws.Range("P4").Value = 0.3083291105451986
$ws.Range("Q4").Value = 0.1821826722906667
$ws.Range("R4").Value = 1.639644050616
$ws.Range("S4").Value = 0.02459402256146521
$ws.Range("T4").Value = 0.0245940225614652

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Egf"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3966103333333333
$ws.Range("H5").Value = 1.189831
$ws.Range("I5").Value = 0.4967728144006086
$ws.Range("J5").Value = 0.4967728144006086
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.179771666666667
$ws.Range("N5").Value = 6.539315
$ws.Range("O5").Value = 0.2349306639444428
$ws.Range("P5").Value = 0.2349306639444428
$ws.Range("Q5").Value = 0.8645199673072224
$ws.Range("R5").Value = 7.780679705765001
$ws.Range("S5").Value = 0.1167071671166844
$ws.Range("T5").Value = 0.1167071671166844

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Egf"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3966103333333333
$ws.Range("H6").Value = 1.189831
$ws.Range("I6").Value = 0.4967728144006086
$ws.Range("J6").Value = 0.4967728144006086
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.237801
$ws.Range("N6").Value = 12.713403
$ws.Range("O6").Value = 0.4567402255103586
$ws.Range("P6").Value = 0.4567402255103586
$ws.Range("Q6").Value = 1.680755667210333
$ws.Range("R6").Value = 15.126801004893
$ws.Range("S6").Value = 0.2268961272767495
$ws.Range("T6").Value = 0.2268961272767495

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Egf"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3966103333333333
$ws.Range("H7").Value = 1.189831
$ws.Range("I7").Value = 0.4967728144006086
$ws.Range("J7").Value = 0.4967728144006086
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.860789
$ws.Range("N7").Value = 8.582367000000001
$ws.Range("O7").Value = 0.3083291105451987
$ws.Range("P7").Value = 0.3083291105451986
$ws.Range("Q7").Value = 1.134618478886334
$ws.Range("R7").Value = 10.211566309977
$ws.Range("S7").Value = 0.1531695200071747
$ws.Range("T7").Value = 0.1531695200071747

$ws.Range("A8").Value = "ECs"
$ws.Range("B8").Value = "Egf"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3380806666666666
$ws.Range("H8").Value = 1.014242
$ws.Range("I8").Value = 0.4234616956721602
$ws.Range("J8").Value = 0.4234616956721601
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 2.179771666666667
$ws.Range("N8").Value = 6.539315
$ws.Range("O8").Value = 0.2349306639444428
$ws.Range("P8").Value = 0.2349306639444428
$ws.Range("Q8").Value = 0.7369386582477778
$ws.Range("R8").Value = 6.632447924229999
$ws.Range("S8").Value = 0.09948413731930017
$ws.Range("T8").Value = 0.09948413731930014

$ws.Range("A9").Value = "ECs"
$ws.Range("B9").Value = "Egf"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3380806666666666
$ws.Range("H9").Value = 1.014242
$ws.Range("I9").Value = 0.4234616956721602
$ws.Range("J9").Value = 0.4234616956721601
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.237801
$ws.Range("N9").Value = 12.713403
$ws.Range("O9").Value = 0.4567402255103586
$ws.Range("P9").Value = 0.4567402255103586
$ws.Range("Q9").Value = 1.432718587280667
$ws.Range("R9").Value = 12.894467285526
$ws.Range("S9").Value = 0.1934119903763013
$ws.Range("T9").Value = 0.1934119903763012

$ws.Range("A10").Value = "ECs"
$ws.Range("B10").Value = "Egf"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3380806666666666
$ws.Range("H10").Value = 1.014242
$ws.Range("I10").Value = 0.4234616956721602
$ws.Range("J10").Value = 0.4234616956721601
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.860789
$ws.Range("N10").Value = 8.582367000000001
$ws.Range("O10").Value = 0.3083291105451987
$ws.Range("P10").Value = 0.3083291105451986
$ws.Range("Q10").Value = 0.9671774523126667
$ws.Range("R10").Value = 8.704597070814
$ws.Range("S10").Value = 0.1305655679765588
$ws.Range("T10").Value = 0.1305655679765587
